# Add two new rows of task data to the worksheet (rows 6 and 7),
# matching the existing table layout (A: id, B: name, C: start date,
# D: end date, E: status icon).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 6; A = 6; B = "Prueba nueva"; C = "20-10-2024"; D = "22-12-2024"; E = "❌" },
    @{ Row = 7; A = 7; B = "Prueba";       C = "20-10-2024"; D = "22-11-2024"; E = "❌" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E

    # New cells otherwise pick up an explicit style index inherited from
    # the column defaults; reset to Normal so the cells stay unstyled,
    # matching the rest of the table's cells.
    $ws.Range($ws.Cells.Item($r.Row, 1), $ws.Cells.Item($r.Row, 5)).Style = "Normal"
}
